$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row data swaps/rotations (column A "id" stays put; columns B:AC carry the match data) ---

# Rows 6 <-> 7
$v6  = $ws.Range("B6:AC6").Value()
$v7  = $ws.Range("B7:AC7").Value()
$ws.Range("B6:AC6").Value = $v7
$ws.Range("B7:AC7").Value = $v6

# Rows 30 <-> 31
$v30 = $ws.Range("B30:AC30").Value()
$v31 = $ws.Range("B31:AC31").Value()
$ws.Range("B30:AC30").Value = $v31
$ws.Range("B31:AC31").Value = $v30

# Rows 43, 44, 45 rotate: new43 = old45, new44 = old43, new45 = old44
$v43 = $ws.Range("B43:AC43").Value()
$v44 = $ws.Range("B44:AC44").Value()
$v45 = $ws.Range("B45:AC45").Value()
$ws.Range("B43:AC43").Value = $v45
$ws.Range("B44:AC44").Value = $v43
$ws.Range("B45:AC45").Value = $v44

# Rows 59, 60, 61 rotate: new59 = old61, new60 = old59, new61 = old60
$v59 = $ws.Range("B59:AC59").Value()
$v60 = $ws.Range("B60:AC60").Value()
$v61 = $ws.Range("B61:AC61").Value()
$ws.Range("B59:AC59").Value = $v61
$ws.Range("B60:AC60").Value = $v59
$ws.Range("B61:AC61").Value = $v60

# Rows 65 <-> 66
$v65 = $ws.Range("B65:AC65").Value()
$v66 = $ws.Range("B66:AC66").Value()
$ws.Range("B65:AC65").Value = $v66
$ws.Range("B66:AC66").Value = $v65

# Rows 95 <-> 96
$v95 = $ws.Range("B95:AC95").Value()
$v96 = $ws.Range("B96:AC96").Value()
$ws.Range("B95:AC95").Value = $v96
$ws.Range("B96:AC96").Value = $v95

# Rows 109 <-> 110
$v109 = $ws.Range("B109:AC109").Value()
$v110 = $ws.Range("B110:AC110").Value()
$ws.Range("B109:AC109").Value = $v110
$ws.Range("B110:AC110").Value = $v109

# --- Remove the last 4 (future/unplayed) fixtures, rows 135-138 ---
$ws.Range("A135:AC138").Delete()
